$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: update only the changed odds cells (per diff) ---
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 2.63
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.93
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 15
$ws.Range("AD3").Value = 7
$ws.Range("AL3").Value = 29
$ws.Range("AQ3").Value = 34
$ws.Range("AZ3").Value = 67

# --- Row 4: replaced in full with new match data (old row 5 shifted up, with tweaks) ---
$ws.Range("A4").Value = "hdhSltk1"
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "06/11/2024"
$c.Style = "Normal"
$ws.Range("C4").Value = "13:00"
$ws.Range("D4").Value = "LITHUANIA - A LYGA"
$ws.Range("E4").Value = "Hegelmann"
$ws.Range("F4").Value = "Transinvest"
$ws.Range("G4").Value = 2.55
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 2.32
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 2.18
$ws.Range("L4").Value = 2.9
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.19
$ws.Range("P4").Value = 3.62
$ws.Range("Q4").Value = 1.72
$ws.Range("R4").Value = 1.9
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 2.52
$ws.Range("U4").Value = 1.6
$ws.Range("V4").Value = 2.27
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 11.25
$ws.Range("Y4").Value = 8.25
$ws.Range("Z4").Value = 22
$ws.Range("AA4").Value = 16
$ws.Range("AB4").Value = 21
$ws.Range("AC4").Value = 11.75
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 11
$ws.Range("AF4").Value = 40
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 7.8
$ws.Range("AI4").Value = 10.25
$ws.Range("AJ4").Value = 7.9
$ws.Range("AK4").Value = 19
$ws.Range("AL4").Value = 15
$ws.Range("AM4").Value = 21
$ws.Range("AN4").Value = 4.55
$ws.Range("AO4").Value = 13.5
$ws.Range("AP4").Value = 20
$ws.Range("AQ4").Value = 55
$ws.Range("AR4").Value = 90
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.65
$ws.Range("AU4").Value = 6.9
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 4.35
$ws.Range("AX4").Value = 12
$ws.Range("AY4").Value = 19.5
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 80
$ws.Range("BB4").Value = 250

# --- Row 5: delete entirely (its data became row 4) ---
$ws.Rows(5).Delete()
